$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "Nicolas Andres Bazan Antinao"
$find.Replacement.Text = "NICOLAS ANDRES BAZAN ANTINAO"
$find.Forward = $true
$find.Wrap = 1
$find.Format = $false
$find.MatchCase = $false
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null
